$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; this pushes the existing rows 12-43
# down to 13-44 (values, not just blank rows) and grows the used range
# to A1:R44 automatically.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 45177
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = "Alcachofa"
$ws.Range("H12").Value = "Argentina(o)"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = "$/caja 50 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 300
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = "Hortaliza"

# Match the date-number formatting used by the rest of column D.
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
